# Apply edits described by the commit "save figures and results for summary figures."
#
# Summary of changes:
#  - On the "cell_types" sheet, insert a new row for "PBMC" / "#999999"
#    right before the "backup1" / "backup2" rows (which shifts them down
#    by one row), and re-shuffle a couple of existing color values.
#  - Make "cell_types" the active sheet/tab (instead of "validation"),
#    with the cell selection on B4.

$wb = $excel.ActiveWorkbook

$wsCellTypes  = $wb.Worksheets.Item("cell_types")
$wsValidation = $wb.Worksheets.Item("validation")

# Insert a new row above row 14 (which currently holds "backup1"), shifting
# "backup1"/"backup2" down to rows 15/16.
$wsCellTypes.Rows.Item(14).Insert()

# Update the color value for "NK cell" (row 4, column B).
$wsCellTypes.Range("B4").Value = "#ffed6f"

# Update the color value for "other cell" (row 13, column B).
$wsCellTypes.Range("B13").Value = "#bebada"

# Fill in the newly inserted row 14 with the PBMC entry.
$wsCellTypes.Range("A14").Value = "PBMC"
$wsCellTypes.Range("B14").Value = "#999999"

# Restore values for the rows that were pushed down (backup1 / backup2),
# matching their pre-insert contents.
$wsCellTypes.Range("A15").Value = "backup1"
$wsCellTypes.Range("B15").Value = "#ffed6f"
$wsCellTypes.Range("A16").Value = "backup2"
$wsCellTypes.Range("B16").Value = "#ccebc5"

# Make "cell_types" the active sheet with the given selection, and clear
# the previously active tab flag on "validation".
$wsCellTypes.Select()
$wsCellTypes.Range("B4").Select()
$wsValidation.Range("A10").Select()

$wsCellTypes.Activate()
